# "Random with different groups" - the randomizer re-rolled the group
# assignments in column C; row 8 moved from group 4 to group 3.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C8").Value = 3

# Move/restore the active selection to C11 (matches the saved cursor
# position recorded in the workbook after the edit).
$ws.Range("C11").Select() | Out-Null
